$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.321.63"
$ws.Range("E2").Value = "  +3.43%  "
$ws.Range("D3").Value = "3.492.95"
$ws.Range("E3").Value = "  +2.94%  "
$ws.Range("E4").Value = "  -0.03%  "
$c = $ws.Range("D5")
$c.Value = "'581.12"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.62%  "
$c = $ws.Range("D6")
$c.Value = "'163.32"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +5.12%  "
$c = $ws.Range("D7")
$c.Value = "'0.612"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +12.55%  "
$c = $ws.Range("D8")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "3.493.39"
$ws.Range("E9").Value = "  +2.97%  "
$c = $ws.Range("D10")
$c.Value = "'7.26"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.98%  "
$ws.Range("E11").Value = "  +3.70%  "
$c = $ws.Range("D12")
$c.Value = "'0.448"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +3.88%  "
$ws.Range("D13").Value = "4.092.48"
$ws.Range("E13").Value = "  +2.81%  "
$c = $ws.Range("D14")
$c.Value = "'0.134"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("E15").Value = "  +1.71%  "
$c = $ws.Range("D16")
$c.Value = "'28.85"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +6.30%  "
$ws.Range("D17").Value = "65.282.23"
$ws.Range("E17").Value = "  +3.26%  "
$ws.Range("D18").Value = "3.507.97"
$ws.Range("E18").Value = "  +4.55%  "
$ws.Range("E19").Value = "  +3.77%  "
$ws.Range("E20").Value = "  +2.72%  "
$c = $ws.Range("D21")
$c.Value = "'385.73"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.14%  "
$ws.Range("E22").Value = "  +2.70%  "
$c = $ws.Range("D23")
$c.Value = "'0.555"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +5.00%  "
$c = $ws.Range("D24")
$c.Value = "'72.67"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("E26").Value = "  +2.25%  "
$c = $ws.Range("D27")
$c.Value = "'10.08"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +7.02%  "
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +13.02%  "
$c = $ws.Range("D31")
$c.Value = "'6.18"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +2.04%  "
$ws.Range("E32").Value = "  +3.44%  "
$c = $ws.Range("D33")
$c.Value = "'23.79"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.87%  "
$c = $ws.Range("D34")
$c.Value = "'7.22"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +6.44%  "
$ws.Range("E35").Value = "  +12.78%  "
$c = $ws.Range("D36")
$c.Value = "'162.40"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.57%  "
$ws.Range("E37").Value = "  +6.04%  "
$ws.Range("D38").Value = "3.019.27"
$ws.Range("E38").Value = "  +2.27%  "
$c = $ws.Range("D39")
$c.Value = "'0.0785"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +4.26%  "
$c = $ws.Range("D40")
$c.Value = "'27.11"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.33%  "
$c = $ws.Range("D41")
$c.Value = "'6.84"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +7.99%  "
$c = $ws.Range("D42")
$c.Value = "'4.59"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +6.45%  "
$c = $ws.Range("D43")
$c.Value = "'0.0323"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.90%  "
$c = $ws.Range("D44")
$c.Value = "'43.02"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +3.40%  "
$c = $ws.Range("D45")
$c.Value = "'0.782"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +3.05%  "
$c = $ws.Range("D46")
$c.Value = "'25.99"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +12.02%  "
$c = $ws.Range("D47")
$c.Value = "'1.12"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +4.84%  "
$c = $ws.Range("D48")
$c.Value = "'321.61"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +10.41%  "
$c = $ws.Range("D49")
$c.Value = "'6.76"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +6.68%  "
$ws.Range("E51").Value = "  +6.69%  "
